$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2,4)
$c.NumberFormat = "@"
$c.Value = "27.650.90"
$ws.Cells.Item(2,5).Value = "  -0.29%  "

# Row 3
$c = $ws.Cells.Item(3,4)
$c.NumberFormat = "@"
$c.Value = "1.865.02"
$ws.Cells.Item(3,5).Value = "  -0.81%  "

# Row 4
$c = $ws.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = "1.009"
$ws.Cells.Item(4,5).Value = "  +0.06%  "

# Row 5
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "333.65"
$ws.Cells.Item(5,5).Value = "  +0.09%  "

# Row 6
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "1.009"
$ws.Cells.Item(6,5).Value = "  +0.21%  "

# Row 7
$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = "0.4700"
$ws.Cells.Item(7,5).Value = "  -0.03%  "

# Row 8
$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = "0.3922"
$ws.Cells.Item(8,5).Value = "  -0.43%  "

# Row 9
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = "45.84"
$ws.Cells.Item(9,5).Value = "  -3.90%  "

# Row 10
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = "0.07994"
$ws.Cells.Item(10,5).Value = "  -0.82%  "

# Row 11
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "1.004"
$ws.Cells.Item(11,5).Value = "  -2.37%  "

# Row 12
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = "21.81"
$ws.Cells.Item(12,5).Value = "  -1.44%  "

# Row 13
$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = "1.874.42"
$ws.Cells.Item(13,5).Value = "  -0.53%  "

# Row 14
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = "5.998"
$ws.Cells.Item(14,5).Value = "  +0.34%  "

# Row 15
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = "7.252"
$ws.Cells.Item(15,5).Value = "  +1.71%  "

# Row 16
$c = $ws.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value = "1.011"
$ws.Cells.Item(16,5).Value = "  +0.14%  "

# Row 17
$ws.Cells.Item(17,5).Value = "  +1.43%  "

# Row 18
$c = $ws.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value = "0.06728"
$ws.Cells.Item(18,5).Value = "  +0.54%  "

# Row 19
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "0.00001043"
$ws.Cells.Item(19,5).Value = "  -0.60%  "

# Row 20
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = "17.10"
$ws.Cells.Item(20,5).Value = "  -1.28%  "

# Row 21
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = "1.009"
$ws.Cells.Item(21,5).Value = "  +0.23%  "

# Row 22
$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = "27.657.46"
$ws.Cells.Item(22,5).Value = "  -0.27%  "

# Row 23
$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value = "5.461"
$ws.Cells.Item(23,5).Value = "  -1.45%  "

# Row 24
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = "10.90"
$ws.Cells.Item(24,5).Value = "  -1.16%  "

# Row 25
$c = $ws.Cells.Item(25,4)
$c.NumberFormat = "@"
$c.Value = "2.308"
$ws.Cells.Item(25,5).Value = "  -0.42%  "

# Row 26
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = "2.097.76"
$ws.Cells.Item(26,5).Value = "  -0.48%  "

# Row 27
$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = "159.17"
$ws.Cells.Item(27,5).Value = "  -0.69%  "

# Row 28
$c = $ws.Cells.Item(28,4)
$c.NumberFormat = "@"
$c.Value = "19.78"
$ws.Cells.Item(28,5).Value = "  -2.24%  "

# Row 29
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = "2.152"
$ws.Cells.Item(29,5).Value = "  +2.27%  "

# Row 30
$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value = "5.438"
$ws.Cells.Item(30,5).Value = "  -2.64%  "

# Row 31
$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value = "121.80"
$ws.Cells.Item(31,5).Value = "  -0.11%  "

# Row 32
$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = "0.9808"
$ws.Cells.Item(32,5).Value = "  -0.45%  "

# Row 33
$c = $ws.Cells.Item(33,4)
$c.NumberFormat = "@"
$c.Value = "0.09495"
$ws.Cells.Item(33,5).Value = "  +0.10%  "

# Row 34
$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = "3.610"
$ws.Cells.Item(34,5).Value = "  -0.15%  "

# Row 36
$ws.Cells.Item(36,5).Value = "  -7.91%  "

# Row 37
$c = $ws.Cells.Item(37,4)
$c.NumberFormat = "@"
$c.Value = "0.06058"
$ws.Cells.Item(37,5).Value = "  -1.37%  "

# Row 38
$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value = "0.02231"
$ws.Cells.Item(38,5).Value = "  -1.46%  "

# Row 39
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = "8.325"
$ws.Cells.Item(39,5).Value = "  +2.63%  "

# Row 40
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "1.194"
$ws.Cells.Item(40,5).Value = "  -2.75%  "

# Row 41
$ws.Cells.Item(41,5).Value = "  +0.33%  "

# Row 42
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = "0.5961"
$ws.Cells.Item(42,5).Value = "  -0.60%  "

# Row 43
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = "0.1884"
$ws.Cells.Item(43,5).Value = "  -0.92%  "

# Row 44
$ws.Cells.Item(44,5).Value = "  +0.00%  "

# Row 45
$c = $ws.Cells.Item(45,4)
$c.NumberFormat = "@"
$c.Value = "1.250"
$ws.Cells.Item(45,5).Value = "  -0.76%  "

# Row 46
$ws.Cells.Item(46,5).Value = "  -1.42%  "

# Row 47
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = "12.18"
$ws.Cells.Item(47,5).Value = "  -0.21%  "

# Row 48
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = "1.924"
$ws.Cells.Item(48,5).Value = "  -1.23%  "

# Row 49
$ws.Cells.Item(49,2).Value = "Cronos"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Cells.Item(49,4)
$c.NumberFormat = "@"
$c.Value = "0.06754"
$ws.Cells.Item(49,5).Value = "  -2.30%  "

# Row 50
$ws.Cells.Item(50,2).Value = "Quant"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value = "112.09"
$ws.Cells.Item(50,5).Value = "  -2.28%  "

# Row 51
$ws.Cells.Item(51,2).Value = "PancakeSwap"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "3.127"
$ws.Cells.Item(51,5).Value = "  -8.12%  "
